# "Add auto type detection and more types"
# Adds a new "types" worksheet (after "basic") that demonstrates auto
# type-detection across a grid of value kinds: string, date, datetime,
# time, bool, integer, float, currency, percentage, scientific.

$wb = $excel.ActiveWorkbook

# --- create the new sheet at the end of the tab strip ---------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "types"

# --- column A : string ------------------------------------------------
$ws.Range("A1").Value = "string"
$ws.Range("A2").Value = 1
$ws.Range("A3").NumberFormat = "mm-dd-yy"
$ws.Range("A3").Value = 45782
$ws.Range("A4").Value = 2.4
$ws.Range("A5").Value = "test"

# --- column B : date ----------------------------------------------------
$ws.Range("B1").Value = "date"
$ws.Range("B2").NumberFormat = "mm-dd-yy"
$ws.Range("B2").Value = 45782
$ws.Range("B3").NumberFormat = "@"
$ws.Range("B3").Value = "2025-05-05"

# --- column E : bool ------------------------------------------------------
$ws.Range("E1").Value = "bool"
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 0
$ws.Range("E4").Value = $true
$ws.Range("E5").Value = "f"

# --- column F : integer -----------------------------------------------
$ws.Range("F1").Value = "integer"
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 2
$ws.Range("F4").Value = 4

# --- column G : float -------------------------------------------------
$ws.Range("G1").Value = "float"
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 4.5

# --- column C : datetime ------------------------------------------------
$ws.Range("C1").Value = "datetime"
$ws.Range("C2").NumberFormat = "mm-dd-yy"
$ws.Range("C2").Value = 45782.625
$ws.Range("C3").NumberFormat = "mm-dd-yy"
$ws.Range("C3").Value = 45782
$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "2025-05-05"

# --- column H : currency ------------------------------------------------
$ws.Range("H1").Value = "currency"
$ws.Range("H2").NumberFormat = "$#,##0.00"
$ws.Range("H2").Value = 5
$ws.Range("H3").Value = 10

# --- column I : percentage ----------------------------------------------
$ws.Range("I1").Value = "percentage"
$ws.Range("I2").NumberFormat = "0.00%"
$ws.Range("I2").Value = 0.2
$ws.Range("I3").Value = 0.1
$ws.Range("I4").Value = 50

# --- column J : scientific ----------------------------------------------
$ws.Range("J1").Value = "scientific"
$ws.Range("J2").NumberFormat = "0.00E+00"
$ws.Range("J2").Value = 10
$ws.Range("J3").NumberFormat = "0.00E+00"
$ws.Range("J3").Value = 10000
$ws.Range("J4").Value = 200
$ws.Range("J5").Value = 0.2

# --- column D : time ------------------------------------------------------
$ws.Range("D1").Value = "time"
$ws.Range("D2").NumberFormat = "h:mm"
$ws.Range("D2").Value = 0.66666666666666663
$ws.Range("D3").NumberFormat = "mm-dd-yy"
$ws.Range("D3").Value = "17h00"
$ws.Range("D4").NumberFormat = "@"

# --- column widths (auto-fit look from the original authoring) -----------
$ws.Columns.Item(1).ColumnWidth = 23.59
$ws.Columns.Item(2).ColumnWidth = 17.88
$ws.Columns.Item(3).ColumnWidth = 12.59
$ws.Columns.Item(4).ColumnWidth = 12.59
$ws.Columns.Item(5).ColumnWidth = 13.45
$ws.Columns.Item(6).ColumnWidth = 16.45
$ws.Columns.Item(8).ColumnWidth = 23.74
$ws.Columns.Item(9).ColumnWidth = 19.45
$ws.Columns.Item(10).ColumnWidth = 17.45

# --- make the new sheet the active / selected tab -------------------------
$ws.Activate() | Out-Null
$ws.Range("J6").Select() | Out-Null
